# Update gh-pages output: refresh "想去人数" counters for several listed
# events, split the former "南宁·蔚蓝档案only" row into two rows (a new
# "南宁·火影忍者only" entry followed by the updated "蔚蓝档案only" entry)
# on both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

function Update-ExpoSheet {
    # positional params only:
    #   ws, rAp, rStar, rHuaYuan, rAb, rLiangYa, rBlue, rNew
    param($ws, $rAp, $rStar, $rHuaYuan, $rAb, $rLiangYa, $rBlue, $rNew)

    # --- simple counter refreshes (F column = "想去人数") ---
    $ws.Cells.Item($rAp, 6).Value = 4428
    $ws.Cells.Item($rStar, 6).Value = 57
    $ws.Cells.Item($rHuaYuan, 6).Value = 165
    $ws.Cells.Item($rAb, 6).Value = 1657
    $ws.Cells.Item($rLiangYa, 6).Value = 3570

    # --- insert a fresh row right after the "蔚蓝档案only" row; that new
    #     row will hold a copy of the "蔚蓝档案only" data (with an updated
    #     count), while the original row turns into a brand new
    #     "火影忍者only" entry ---
    $ws.Rows.Item($rNew).Insert()

    # copy the "蔚蓝档案only" row's formatting (border/bold/center on col A)
    # down onto the newly inserted row
    $ws.Cells.Item($rBlue, 1).Copy()
    $ws.Cells.Item($rNew, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # fill the new row with what used to be in the "蔚蓝档案only" row,
    # bumping the "想去人数" count from 235 to 237
    $ws.Cells.Item($rNew, 1).Value = $rNew - 1
    $bCell = $ws.Cells.Item($rNew, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = "2024-08-03"
    $ws.Cells.Item($rNew, 3).Value = "南宁·蔚蓝档案only"
    $ws.Cells.Item($rNew, 4).Value = "亭洪路45号 百益上河城"
    $ws.Cells.Item($rNew, 5).Value = "2024.08.03 09:00-08.03 17:00"
    $ws.Cells.Item($rNew, 6).Value = 237
    $ws.Cells.Item($rNew, 7).Value = 68
    $ws.Cells.Item($rNew, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85370"
    $ws.Cells.Item($rNew, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg"

    # turn the original row into the new "南宁·火影忍者only" entry
    $ws.Cells.Item($rBlue, 3).Value = "南宁·火影忍者only"
    $ws.Cells.Item($rBlue, 4).Value = "厢竹大道65号 桔子酒店"
    $ws.Cells.Item($rBlue, 5).Value = "2024.08.03 10:00-08.03 17:00"
    $ws.Cells.Item($rBlue, 6).Value = 5
    $ws.Cells.Item($rBlue, 7).Value = 68
    $ws.Cells.Item($rBlue, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86994"
    $ws.Cells.Item($rBlue, 9).Value = "//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg"
}

# Sheet "展览": 蔚蓝档案only currently on row 15, new row goes in at 16
$wsExpo = $wb.Worksheets.Item("展览")
Update-ExpoSheet $wsExpo 2 6 11 12 14 15 16

# Sheet "全部类型": 蔚蓝档案only currently on row 19, new row goes in at 20
$wsAll = $wb.Worksheets.Item("全部类型")
Update-ExpoSheet $wsAll 2 7 13 16 18 19 20
